$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.398.87"
$ws.Range("E2").Value = "  +3.98%  "
$ws.Range("D3").Value = "1.587.64"
$ws.Range("E3").Value = "  +1.23%  "
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").Value = "'214.29"
$ws.Range("E5").Value = "  +1.55%  "
$ws.Range("E6").Value = "  +0.89%  "
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("D8").Value = "'23.90"
$ws.Range("E8").Value = "  +7.99%  "
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("D11").Value = "'0.0888"
$ws.Range("E11").Value = "  +2.13%  "
$ws.Range("D12").Value = "1.813.68"
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("D13").Value = "1.593.13"
$ws.Range("E13").Value = "  +1.62%  "
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("E15").Value = "  +2.37%  "
$ws.Range("D16").Value = "28.396.06"
$ws.Range("E16").Value = "  +4.16%  "
$ws.Range("D17").Value = "'63.08"
$ws.Range("E17").Value = "  +1.28%  "
$ws.Range("D18").Value = "'232.32"
$ws.Range("E18").Value = "  +6.87%  "
$ws.Range("E19").Value = "  +0.70%  "
$ws.Range("D20").Value = "'7.49"
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("E21").Value = "  -0.38%  "
$ws.Range("D22").Value = "'4.11"
$ws.Range("E22").Value = "  -0.96%  "
$ws.Range("E23").Value = "  +1.74%  "
$ws.Range("D24").Value = "'1.95"
$ws.Range("E24").Value = "  +0.74%  "
$ws.Range("D25").Value = "'151.98"
$ws.Range("E25").Value = "  -1.00%  "
$ws.Range("E26").Value = "  +1.11%  "
$ws.Range("D27").Value = "'6.62"
$ws.Range("E27").Value = "  -0.46%  "
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("E30").Value = "  +0.44%  "
$ws.Range("D31").Value = "'0.0473"
$ws.Range("E31").Value = "  +0.27%  "
$ws.Range("E32").Value = "  -0.28%  "
$ws.Range("D33").Value = "'3.16"
$ws.Range("E33").Value = "  -0.48%  "
$ws.Range("D34").Value = "1.413.44"
$ws.Range("E34").Value = "  -2.26%  "
$ws.Range("E35").Value = "  -1.69%  "
$ws.Range("D36").Value = "'1.05"
$ws.Range("E36").Value = "  -4.75%  "
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("E38").Value = "  +0.41%  "
$ws.Range("E39").Value = "  +8.19%  "
$ws.Range("D40").Value = "'0.543"
$ws.Range("E40").Value = "  +1.83%  "
$ws.Range("E41").Value = "  +0.77%  "
$ws.Range("D42").Value = "'5.76"
$ws.Range("E42").Value = "  -2.33%  "
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("D44").Value = "'0.980"
$ws.Range("D46").Value = "'64.62"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").Value = "1.726.74"
$ws.Range("E47").Value = "  +1.38%  "
$ws.Range("D48").Value = "'87.58"
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("E49").Value = "  +5.71%  "
$ws.Range("E50").Value = "  -0.76%  "
$ws.Range("D51").Value = "'39.14"
$ws.Range("E51").Value = "  +15.42%  "
